$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.576158285140991
$ws.Range("B1").Value = 1.697693467140198
$ws.Range("C1").Value = 2.133358478546143
$ws.Range("D1").Value = 2.221862077713013
$ws.Range("E1").Value = 1.434656500816345
